# Lock in current version of the cues sequence: refresh the word list and
# the paired image assignments (and, where needed, the matching category
# label) in Sheet1 for rows 2-49. Columns: A=word, B=image, C=category.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(2,  "wenden",    "none"),
    @(3,  "kümmern",   "face/face013.jpg"),
    @(4,  "schulden",  "flower/flower020.jpg"),
    @(5,  "rufen",     "none"),
    @(6,  "schreiben", "flower/flower012.jpg"),
    @(7,  "kleben",    "flower/flower017.jpg"),
    @(8,  "enden",     "none"),
    @(9,  "hoffen",    "face/face002.jpg"),
    @(10, "wagen",     "face/face025.jpg"),
    @(11, "opfern",    "none"),
    @(12, "danken",    "flower/flower014.jpg"),
    @(13, "handeln",   "flower/flower005.jpg"),
    @(14, "schalten",  "none"),
    @(15, "machen",    "flower/flower031.jpg"),
    @(16, "schütteln", "face/face011.jpg"),
    @(17, "klagen",    "none"),
    @(18, "bilden",    "flower/flower019.jpg"),
    @(19, "mauern",    "face/face020.jpg"),
    @(20, "sparen",    "none"),
    @(21, "tollen",    "flower/flower013.jpg"),
    @(22, "schicken",  "flower/flower003.jpg"),
    @(23, "hören",     "none"),
    @(24, "schleppen", "face/face004.jpg"),
    @(25, "passen",    "flower/flower006.jpg"),
    @(26, "drohen",    "none"),
    @(27, "regnen",    "flower/flower023.jpg"),
    @(28, "töten",     "face/face017.jpg"),
    @(29, "orten",     "none"),
    @(30, "münzen",    "flower/flower015.jpg"),
    @(31, "bauen",     "face/face018.jpg"),
    @(32, "weigern",   "none"),
    @(33, "deuten",    "flower/flower000.jpg"),
    @(34, "platzen",   "face/face003.jpg"),
    @(35, "ächzen",    "none"),
    @(36, "gelten",    "flower/flower026.jpg"),
    @(37, "rühren",    "face/face000.jpg"),
    @(38, "kosten",    "none"),
    @(39, "seufzen",   "flower/flower010.jpg"),
    @(40, "leuchten",  "face/face007.jpg"),
    @(41, "dauern",    "none"),
    @(42, "lächeln",   "face/face012.jpg"),
    @(43, "testen",    "face/face027.jpg"),
    @(44, "stören",    "none"),
    @(45, "küssen",    "face/face029.jpg"),
    @(46, "fühlen",    "flower/flower024.jpg"),
    @(47, "bremsen",   "none"),
    @(48, "decken",    "face/face031.jpg"),
    @(49, "stopfen",   "face/face014.jpg")
)

foreach ($entry in $data) {
    $row = $entry[0]
    $word = $entry[1]
    $image = $entry[2]

    if ($image -eq "none") {
        $category = "none"
    } else {
        $category = $image.Split("/")[0]
    }

    $ws.Cells.Item($row, 1).Value = $word
    $ws.Cells.Item($row, 2).Value = $image
    $ws.Cells.Item($row, 3).Value = $category
}
